$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gasto Capital")

# Insert two new rows before the current "2.- Servicios" row (row 12),
# shifting the existing "2.- Servicios" / "MANTENIMIENTO DE EXTINTORES"
# rows down to rows 14/15.
$ws.Rows("12:13").Insert()

# New row 12: first filler material line item
$ws.Range("B12").Value = "FILLER DE 0.0040 A 0.009 25 HOJAS X 4"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 58.98

# New row 13: second filler material line item
$ws.Range("B13").Value = "FILLER DE 0.050 A 1.000  20 HOJAS X 12"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 81.42
